# Update molar volume calculations: refresh UMAP Component 1 (AE) and
# UMAP Component 2 (AF) values for rows 2-70 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aeValues = @(3.048986434936523, 3.158958435058594, 3.253911733627319, 3.339487314224243, 3.268934011459351, 2.643489360809326, 2.766763687133789, 2.947901725769043, 3.03215479850769, 3.107831478118896, 3.298761129379272, 3.354927778244019, 1.762137293815613, 2.164197683334351, 2.394794702529907, 2.70711612701416, 2.857454061508179, 3.000916719436646, 3.196066379547119, 3.363282442092896, 3.431832075119019, 1.478200554847717, 1.718989968299866, 2.011849880218506, 2.246812582015991, 2.483688354492188, 2.706009864807129, 3.023434162139893, 3.203835248947144, 3.428520202636719, 1.253238201141357, 1.428144097328186, 1.646275401115417, 1.843245983123779, 2.163870811462402, 2.434584617614746, 2.715826511383057, 2.960421323776245, 3.156507253646851, 1.090214490890503, 1.276663303375244, 1.294190764427185, 1.58684766292572, 1.792916893959045, 2.077682971954346, 2.335130214691162, 2.698145389556885, 3.068830251693726, 1.119083642959595, 0.9675148129463196, 1.146602749824524, 1.300629377365112, 1.578214526176453, 1.819749593734741, 2.028735399246216, 2.44606876373291, 2.717647552490234, 1.166656136512756, 0.9909321665763855, 1.040380954742432, 1.263261318206787, 1.546685695648193, 1.877916216850281, 2.220521688461304, 1.139411330223083, 1.068029403686523, 1.139615297317505, 1.44741427898407, 1.525559544563293)
$afValues = @(14.19855213165283, 14.67528915405273, 14.78826999664307, 15.0087947845459, 15.4818696975708, 14.26635646820068, 14.6733865737915, 14.93699264526367, 15.22404766082764, 15.38049697875977, 15.74327278137207, 16.12627220153809, 14.24727058410645, 14.80567073822021, 15.32405567169189, 15.22656726837158, 15.39282512664795, 15.81495380401611, 16.19318008422852, 16.40904808044434, 16.70996475219727, 14.68713569641113, 14.94417858123779, 15.29533672332764, 15.39813709259033, 15.819580078125, 16.15385818481445, 16.43810844421387, 16.62729835510254, 16.89436721801758, 14.86218929290771, 15.21390724182129, 15.43618869781494, 15.77967262268066, 16.05665588378906, 16.35312843322754, 16.59934234619141, 16.86336326599121, 17.02183532714844, 15.0966968536377, 15.32936000823975, 15.5907564163208, 15.97738075256348, 16.12858200073242, 16.52140235900879, 16.83902740478516, 16.9146785736084, 17.04213523864746, 15.52112483978271, 15.61314487457275, 15.75729751586914, 16.19294166564941, 16.5205135345459, 16.58085632324219, 16.90520286560059, 16.96755027770996, 17.05212020874023, 15.84199905395508, 16.08383369445801, 16.37764549255371, 16.6052360534668, 16.79434394836426, 16.81621551513672, 16.89978408813477, 16.07024192810059, 16.1440486907959, 16.52154922485352, 16.44603729248047, 16.81131744384766)

$startRow = 2
for ($i = 0; $i -lt $aeValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 31).Value = $aeValues[$i]   # column AE = 31
    $ws.Cells.Item($row, 32).Value = $afValues[$i]   # column AF = 32
}

Write-Output "Updated AE2:AF70 with new molar volume derived UMAP values"
